$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new blank row at position 18 (old blank gap row 18 and the two
# trailing footnote rows all shift down by one, matching the target layout:
# old row 18 (blank) -> 19, old row 19 -> 20, old row 20 -> 21).
$ws.Rows("18:18").Insert()

# Fill the newly inserted row with the "Schlusspräsentation" milestone.
# Column B is written first so the shared-string table appends
# "Schlusspräsentation" before "?" (matches the target sharedStrings order).
$ws.Cells.Item(18, 2).Value = "Schlusspräsentation"
$ws.Cells.Item(18, 2).Font.Bold = $true

$ws.Cells.Item(18, 1).Value = "?"
$ws.Cells.Item(18, 1).Font.Bold = $true
$ws.Cells.Item(18, 1).NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"
$ws.Cells.Item(18, 1).HorizontalAlignment = -4152

# Restore the active selection to B15, as in the target workbook.
$ws.Range("B15").Select()
